# Updated cryptos list with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values
# for rows 2-51 of Sheet1, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D:E cells are treated as plain text so that values such as
# "62.434.05" or "  -3.03%  " are not reinterpreted by Excel as dates,
# numbers, or trimmed of their surrounding whitespace.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "62.434.05"
$ws.Cells.Item(2, 5).Value = "  -3.03%  "
$ws.Cells.Item(3, 4).Value = "3.173.27"
$ws.Cells.Item(3, 5).Value = "  -5.09%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "523.41"
$ws.Cells.Item(5, 5).Value = "  +0.03%  "
$ws.Cells.Item(6, 4).Value = "170.48"
$ws.Cells.Item(6, 5).Value = "  -7.56%  "
$ws.Cells.Item(7, 4).Value = "0.593"
$ws.Cells.Item(7, 5).Value = "  +0.39%  "
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.07%  "
$ws.Cells.Item(9, 4).Value = "3.174.75"
$ws.Cells.Item(9, 5).Value = "  -4.90%  "
$ws.Cells.Item(10, 4).Value = "0.601"
$ws.Cells.Item(10, 5).Value = "  -2.20%  "
$ws.Cells.Item(11, 4).Value = "52.82"
$ws.Cells.Item(11, 5).Value = "  -7.44%  "
$ws.Cells.Item(12, 4).Value = "0.131"
$ws.Cells.Item(12, 5).Value = "  +0.36%  "
$ws.Cells.Item(13, 4).Value = "0.0000250"
$ws.Cells.Item(13, 5).Value = "  +0.63%  "
$ws.Cells.Item(14, 4).Value = "9.01"
$ws.Cells.Item(14, 5).Value = "  -1.24%  "
$ws.Cells.Item(15, 4).Value = "3.667.92"
$ws.Cells.Item(15, 5).Value = "  -5.24%  "
$ws.Cells.Item(16, 4).Value = "0.117"
$ws.Cells.Item(16, 5).Value = "  -3.46%  "
$ws.Cells.Item(17, 4).Value = "3.159.46"
$ws.Cells.Item(17, 5).Value = "  -5.39%  "
$ws.Cells.Item(18, 4).Value = "17.14"
$ws.Cells.Item(18, 5).Value = "  -0.10%  "
$ws.Cells.Item(19, 4).Value = "62.199.21"
$ws.Cells.Item(19, 5).Value = "  -2.82%  "
$ws.Cells.Item(20, 4).Value = "10.95"
$ws.Cells.Item(20, 5).Value = "  +0.85%  "
$ws.Cells.Item(21, 4).Value = "0.964"
$ws.Cells.Item(21, 5).Value = "  +1.32%  "
$ws.Cells.Item(22, 4).Value = "363.80"
$ws.Cells.Item(22, 5).Value = "  -2.00%  "
$ws.Cells.Item(23, 4).Value = "11.20"
$ws.Cells.Item(23, 5).Value = "  +5.66%  "
$ws.Cells.Item(24, 4).Value = "3.72"
$ws.Cells.Item(24, 5).Value = "  +1.76%  "
$ws.Cells.Item(25, 4).Value = "80.65"
$ws.Cells.Item(25, 5).Value = "  +1.07%  "
$ws.Cells.Item(26, 4).Value = "3.90"
$ws.Cells.Item(26, 5).Value = "  +5.19%  "
$ws.Cells.Item(27, 4).Value = "6.12"
$ws.Cells.Item(27, 5).Value = "  +3.11%  "
$ws.Cells.Item(28, 4).Value = "2.62"
$ws.Cells.Item(28, 5).Value = "  +0.54%  "
$ws.Cells.Item(29, 4).Value = "11.26"
$ws.Cells.Item(29, 5).Value = "  +0.86%  "
$ws.Cells.Item(30, 4).Value = "8.13"
$ws.Cells.Item(30, 5).Value = "  -1.72%  "
$ws.Cells.Item(31, 4).Value = "636.87"
$ws.Cells.Item(31, 5).Value = "  -2.62%  "
$ws.Cells.Item(32, 4).Value = "28.17"
$ws.Cells.Item(32, 5).Value = "  -1.62%  "
$ws.Cells.Item(33, 4).Value = "6.41"
$ws.Cells.Item(33, 5).Value = "  -3.67%  "
$ws.Cells.Item(34, 4).Value = "11.27"
$ws.Cells.Item(34, 5).Value = "  +2.64%  "
$ws.Cells.Item(35, 4).Value = "0.105"
$ws.Cells.Item(35, 5).Value = "  +2.44%  "
$ws.Cells.Item(36, 4).Value = "56.45"
$ws.Cells.Item(36, 5).Value = "  -3.78%  "
$ws.Cells.Item(37, 5).Value = "  -0.04%  "
$ws.Cells.Item(38, 4).Value = "36.89"
$ws.Cells.Item(38, 5).Value = "  +3.19%  "
$ws.Cells.Item(39, 4).Value = "0.372"
$ws.Cells.Item(39, 5).Value = "  +0.53%  "
$ws.Cells.Item(40, 5).Value = "  +0.28%  "
$ws.Cells.Item(41, 4).Value = "0.0₃0702"
$ws.Cells.Item(41, 5).Value = "  +15.30%  "
$ws.Cells.Item(42, 5).Value = "  +0.06%  "
$ws.Cells.Item(43, 4).Value = "2.873.98"
$ws.Cells.Item(43, 5).Value = "  +4.36%  "
$ws.Cells.Item(44, 4).Value = "2.51"
$ws.Cells.Item(44, 5).Value = "  +10.36%  "
$ws.Cells.Item(45, 4).Value = "2.90"
$ws.Cells.Item(45, 5).Value = "  +13.98%  "
$ws.Cells.Item(46, 5).Value = "  +2.43%  "
$ws.Cells.Item(47, 4).Value = "0.0391"
$ws.Cells.Item(47, 5).Value = "  +3.24%  "
$ws.Cells.Item(48, 4).Value = "2.56"
$ws.Cells.Item(48, 5).Value = "  -5.07%  "
$ws.Cells.Item(49, 5).Value = "  +6.79%  "
$ws.Cells.Item(50, 4).Value = "0.123"
$ws.Cells.Item(50, 5).Value = "  -0.42%  "
$ws.Cells.Item(51, 4).Value = "134.09"
$ws.Cells.Item(51, 5).Value = "  +0.07%  "
